$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert two new columns before D (shifts old D:K -> F:M)
$ws.Columns("D:E").Insert()

# Apply number-cell formatting (Verdana 12, #,##0, right-aligned) to full data block
$dataRange = $ws.Range("D5:E102")
$dataRange.NumberFormat = "#,##0"
$dataRange.HorizontalAlignment = -4152
$dataRange.Font.Name = "Verdana"
$dataRange.Font.Size = 12
$dataRange.Font.Bold = $false

# Date header rows (7, 38, 80) use bold Verdana + custom date format, general alignment
$dateRows = @(7, 38, 80)
foreach ($r in $dateRows) {
    $rng = $ws.Range("D$($r):E$($r)")
    $rng.NumberFormat = "[$-409]d-mmm-yy;@"
    $rng.Font.Bold = $true
    $rng.HorizontalAlignment = 1
}

# New column D/E values (two additional quarters of data prepended)
$ws.Range("D7").Value = 43465
$ws.Range("E7").Value = 43373
$ws.Range("D8").Value = 207300
$ws.Range("E8").Value = 201500
$ws.Range("D9").Value = 18900
$ws.Range("E9").Value = 29300
$ws.Range("D10").Value = 188400
$ws.Range("E10").Value = 172200
$ws.Range("D12").Value = "NA"
$ws.Range("E12").Value = "NA"
$ws.Range("D13").Value = 0
$ws.Range("E13").Value = 0
$ws.Range("D14").Value = 0
$ws.Range("E14").Value = 0
$ws.Range("D15").Value = 3800
$ws.Range("E15").Value = 3500
$ws.Range("D17").Value = 181500
$ws.Range("E17").Value = 188800
$ws.Range("D18").Value = 25800
$ws.Range("E18").Value = 12700
$ws.Range("D20").Value = -900
$ws.Range("E20").Value = -300
$ws.Range("D21").Value = 28700
$ws.Range("E21").Value = 15900
$ws.Range("D22").Value = 20900
$ws.Range("E22").Value = 19800
$ws.Range("D23").Value = 4000
$ws.Range("E23").Value = -7400
$ws.Range("D24").Value = -100
$ws.Range("E24").Value = -3200
$ws.Range("D25").Value = 0
$ws.Range("E25").Value = 0
$ws.Range("D26").Value = 4100
$ws.Range("E26").Value = -4200
$ws.Range("D27").Value = 4100
$ws.Range("E27").Value = -4200
$ws.Range("D28").Value = 0
$ws.Range("E28").Value = 0
$ws.Range("D29").Value = 100
$ws.Range("E29").Value = "NA"
$ws.Range("D30").Value = 0
$ws.Range("E30").Value = 0
$ws.Range("D31").Value = 0
$ws.Range("E31").Value = 0
$ws.Range("D32").Value = 900
$ws.Range("E32").Value = 300
$ws.Range("D33").Value = 4100
$ws.Range("E33").Value = -4200
$ws.Range("D34").Value = 0
$ws.Range("E34").Value = 0
$ws.Range("D35").Value = 4100
$ws.Range("E35").Value = -4200
$ws.Range("D38").Value = 43465
$ws.Range("E38").Value = 43373
$ws.Range("D41").Value = 58300
$ws.Range("E41").Value = 54800
$ws.Range("D42").Value = 0
$ws.Range("E42").Value = 0
$ws.Range("D43").Value = 37900
$ws.Range("E43").Value = 46400
$ws.Range("D44").Value = 0
$ws.Range("E44").Value = 0
$ws.Range("D45").Value = 11400
$ws.Range("E45").Value = 13200
$ws.Range("D46").Value = 0
$ws.Range("E46").Value = 0
$ws.Range("D47").Value = 561700
$ws.Range("E47").Value = 543000
$ws.Range("D48").Value = 41600
$ws.Range("E48").Value = 36700
$ws.Range("D49").Value = 17700
$ws.Range("E49").Value = 17900
$ws.Range("D50").Value = 0
$ws.Range("E50").Value = 0
$ws.Range("D51").Value = 0
$ws.Range("E51").Value = 0
$ws.Range("D52").Value = 24200
$ws.Range("E52").Value = 23100
$ws.Range("D53").Value = 0
$ws.Range("E53").Value = 0
$ws.Range("D54").Value = 753300
$ws.Range("E54").Value = 736300
$ws.Range("D57").Value = 45000
$ws.Range("E57").Value = 45000
$ws.Range("D58").Value = 0
$ws.Range("E58").Value = 0
$ws.Range("D59").Value = 28900
$ws.Range("E59").Value = 31300
$ws.Range("D60").Value = 0
$ws.Range("E60").Value = 0
$ws.Range("D61").Value = 562600
$ws.Range("E61").Value = 549000
$ws.Range("D62").Value = 0
$ws.Range("E62").Value = 0
$ws.Range("D63").Value = 0
$ws.Range("E63").Value = 0
$ws.Range("D64").Value = 0
$ws.Range("E64").Value = 0
$ws.Range("D65").Value = 0
$ws.Range("E65").Value = 0
$ws.Range("D66").Value = 636500
$ws.Range("E66").Value = 625300
$ws.Range("D68").Value = 0
$ws.Range("E68").Value = 0
$ws.Range("D69").Value = 0
$ws.Range("E69").Value = 0
$ws.Range("D70").Value = 0
$ws.Range("E70").Value = 0
$ws.Range("D71").Value = 0
$ws.Range("E71").Value = 0
$ws.Range("D72").Value = -66500
$ws.Range("E72").Value = -70700
$ws.Range("D73").Value = 0
$ws.Range("E73").Value = 0
$ws.Range("D74").Value = 0
$ws.Range("E74").Value = 0
$ws.Range("D75").Value = 0
$ws.Range("E75").Value = 0
$ws.Range("D76").Value = 116800
$ws.Range("E76").Value = 111100
$ws.Range("D77").Value = 0
$ws.Range("E77").Value = 0
$ws.Range("D80").Value = 43465
$ws.Range("E80").Value = 43373
$ws.Range("D81").Value = 4100
$ws.Range("E81").Value = -4200
$ws.Range("D83").Value = 3800
$ws.Range("E83").Value = 3500
$ws.Range("D84").Value = 0
$ws.Range("E84").Value = 0
$ws.Range("D85").Value = 0
$ws.Range("E85").Value = 0
$ws.Range("D86").Value = 0
$ws.Range("E86").Value = 0
$ws.Range("D87").Value = 0
$ws.Range("E87").Value = 0
$ws.Range("D88").Value = 0
$ws.Range("E88").Value = 0
$ws.Range("D89").Value = 105100
$ws.Range("E89").Value = 94700
$ws.Range("D91").Value = -6100
$ws.Range("E91").Value = -9000
$ws.Range("D92").Value = 0
$ws.Range("E92").Value = 0
$ws.Range("D93").Value = 0
$ws.Range("E93").Value = 0
$ws.Range("D94").Value = -114700
$ws.Range("E94").Value = -132700
$ws.Range("D96").Value = 0
$ws.Range("E96").Value = 0
$ws.Range("D97").Value = 0
$ws.Range("E97").Value = 0
$ws.Range("D98").Value = 0
$ws.Range("E98").Value = 0
$ws.Range("D99").Value = 0
$ws.Range("E99").Value = 0
$ws.Range("D100").Value = 14100
$ws.Range("E100").Value = 23600
$ws.Range("D101").Value = 0
$ws.Range("E101").Value = -100
$ws.Range("D102").Value = 4500
$ws.Range("E102").Value = -14600

# Restated figures in previously-reported quarters (rows 58 and 91)
$ws.Range("F58").Value = 0
$ws.Range("G58").Value = 0
$ws.Range("H58").Value = 0
$ws.Range("I58").Value = 0
$ws.Range("J58").Value = 0
$ws.Range("K58").Value = "NA"
$ws.Range("L58").Value = 0
$ws.Range("M58").Value = 100
$ws.Range("F91").Value = -6400
$ws.Range("G91").Value = -6100
$ws.Range("H91").Value = -4300
$ws.Range("I91").Value = -4100
$ws.Range("J91").Value = -5300
$ws.Range("K91").Value = -6200
$ws.Range("L91").Value = 800
$ws.Range("M91").Value = -13500
